$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the existing "总计" (Total) sheet: prepend a 2022-Q1 summary
#    row and renumber the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 12
$total.Cells.Item(2, 4).Value = 3.78

# Reuse the index column's existing style (bold + centered) instead of
# re-creating it, by copying the format from row 3 (the old row 2).
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$total.Cells.Item(2, 1).Value = 0

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert a brand new "2022-Q1" sheet (fund holdings detail) right
#    after "2021-Q4" and before "总计".
# ---------------------------------------------------------------------
$after = $wb.Worksheets.Item("2021-Q4")
$new = $wb.Worksheets.Add($null, $after)
$new.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $new.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# index, code, name, scale, total stock position, position ratio, held value (100M), position rank
$rows = @(
    @(0,  "516150", "嘉实中证稀土产业ETF",                              "25.17", "99.75", "4.95", "1.2459", 6),
    @(1,  "590002", "中邮核心成长混合",                                  "42.70", "66.61", "1.85", "0.7900", 10),
    @(2,  "516780", "华泰柏瑞中证稀土产业ETF",                          "11.06", "98.70", "4.94", "0.5464", 6),
    @(3,  "001278", "前海开源清洁能源主题精选混合A",                    "8.07",  "68.54", "6.16", "0.4971", 6),
    @(4,  "159715", "易方达中证稀土产业ETF",                            "3.42",  "99.06", "4.90", "0.1676", 6),
    @(5,  "000545", "中邮核心竞争力灵活配置混合",                        "3.43",  "73.83", "4.79", "0.1643", 6),
    @(6,  "159713", "富国中证稀土产业交易型开放式指数证券投资基金",      "3.26",  "99.26", "4.93", "0.1607", 6),
    @(7,  "002360", "前海开源清洁能源主题精选混合C",                    "1.51",  "68.54", "6.16", "0.0930", 6),
    @(8,  "004128", "新疆前海联合泳隆灵活配置混合A",                    "0.86",  "91.05", "5.53", "0.0476", 3),
    @(9,  "007040", "新疆前海联合泳隆灵活配置混合C",                    "0.82",  "91.05", "5.53", "0.0453", 3),
    @(10, "014331", "华泰柏瑞中证稀土产业ETF联接A",                    "0.86",  "24.22", "1.34", "0.0115", 5),
    @(11, "014332", "华泰柏瑞中证稀土产业ETF联接C",                    "0.70",  "24.22", "1.34", "0.0094", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $acell = $new.Cells.Item($rowNum, 1)
    $acell.Value = $data[0]
    $acell.Font.Bold = $true
    $acell.HorizontalAlignment = -4108

    # Columns B..G are stored as text in the source data (keeps leading
    # zeros in fund codes and fixed decimal formatting) - force text type.
    $textRange = $new.Range($new.Cells.Item($rowNum, 2), $new.Cells.Item($rowNum, 7))
    $textRange.NumberFormat = "@"

    $new.Cells.Item($rowNum, 2).Value = $data[1]
    $new.Cells.Item($rowNum, 3).Value = $data[2]
    $new.Cells.Item($rowNum, 4).Value = $data[3]
    $new.Cells.Item($rowNum, 5).Value = $data[4]
    $new.Cells.Item($rowNum, 6).Value = $data[5]
    $new.Cells.Item($rowNum, 7).Value = $data[6]
    $new.Cells.Item($rowNum, 8).Value = $data[7]
}

Write-Output "done"
